$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SOR_Login")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Row 1 / Row 2 swap.
#    Originally: row1 = "Login Page" (A1), row2 = index row 0..11 (A2:L2).
#    Target:     row1 = index row 0..25 (A1:Z1), row2 = "Login Page" (A2).
# ---------------------------------------------------------------------------
$oldA1 = $ws.Range("A1").Value()

# Write the extended index row (0-25) across A1:Z1.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12
$ws.Range("N1").Value = 13
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15
$ws.Range("Q1").Value = 16
$ws.Range("R1").Value = 17
$ws.Range("S1").Value = 18
$ws.Range("T1").Value = 19
$ws.Range("U1").Value = 20
$ws.Range("V1").Value = 21
$ws.Range("W1").Value = 22
$ws.Range("X1").Value = 23
$ws.Range("Y1").Value = 24
$ws.Range("Z1").Value = 25

# Row 2 now holds the "Login Page" text that used to live in A1; it carries no
# special style (unlike the old A2 which had the centred style).
$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = $oldA1

# ---------------------------------------------------------------------------
# 2. New "BC" block (rows 12-13).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "BC_Registration"

$ws.Range("A13").Value = "BC_Name"
$ws.Range("B13").Value = "MaximusInfoNew"
$ws.Range("C13").Value = "PanNo"
$ws.Range("D13").Value = "QXOPC5696Y"
$ws.Range("E13").Value = "AddharNO"
$ws.Range("F13").Value = "'502081483187"
$ws.Range("G13").Value = "AccountNo"
$ws.Range("H13").Value = 8585656999
$ws.Range("I13").Value = "IFSCcode"
$ws.Range("J13").Value = "HDFC0000967"
$ws.Range("K13").Value = "Addess"
$ws.Range("L13").Value = "Thane Maharshtra"
$ws.Range("M13").Value = "PinCode"
$ws.Range("N13").Value = 400001
$ws.Range("O13").Value = "State"
$ws.Range("P13").Value = "Goa"
$ws.Range("Q13").Value = "District"
$ws.Range("R13").Value = "Sangli"
$ws.Range("S13").Value = "City"
$ws.Range("T13").Value = "Satna"
$ws.Range("U13").Value = "EmailID"
$ws.Range("V13").Value = "abc12356@abc.com"
$ws.Hyperlinks.Add($ws.Range("V13"), "mailto:abc12356@abc.com", "", "mailto:abc12356@abc.com", "abc12356@abc.com")
$ws.Range("V13").Style = "Hyperlink"
$ws.Range("W13").Value = "ContactNo"
$ws.Range("X13").Value = 7874563215

# ---------------------------------------------------------------------------
# 3. New "AGG" block (rows 15-16).
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "AGG_Registration"

$ws.Range("A16").Value = "AGG_Name"
$ws.Range("B16").Value = "Maximus"
$ws.Range("C16").Value = "PanNo"
$ws.Range("D16").Value = "QXOPC5696Y"
$ws.Range("E16").Value = "AddharNO"
$ws.Range("F16").Value = "'502081483187"
$ws.Range("G16").Value = "Qualification"
$ws.Range("H16").Value = "12th"
$ws.Range("I16").Value = "AccountNo"
$ws.Range("J16").Value = 8585656999
$ws.Range("K16").Value = "IFSCcode"
$ws.Range("L16").Value = "HDFC0000967"
$ws.Range("M16").Value = "Addess"
$ws.Range("N16").Value = "Thane Maharshtra"
$ws.Range("O16").Value = "PinCode"
$ws.Range("P16").Value = 400001
$ws.Range("Q16").Value = "State"
$ws.Range("R16").Value = "Goa"
$ws.Range("S16").Value = "District"
$ws.Range("T16").Value = "Sangli"
$ws.Range("U16").Value = "City"
$ws.Range("V16").Value = "Satna"
$ws.Range("W16").Value = "EmailID"
$ws.Range("X16").Value = "abc12356@abc.com"
$ws.Hyperlinks.Add($ws.Range("X16"), "mailto:abc12356@abc.com", "", "mailto:abc12356@abc.com", "abc12356@abc.com")
$ws.Range("X16").Style = "Hyperlink"
$ws.Range("Y16").Value = "ContactNo"
$ws.Range("Z16").Value = 7874563215

# ---------------------------------------------------------------------------
# 4. Column width tweaks (values calibrated so the saved OOXML `width` lands
#    on - or as close as this engine's pixel-rounding allows to - the target).
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 25.8
$ws.Columns.Item(8).ColumnWidth = 10.65
$ws.Columns.Item(9).ColumnWidth = 10.65
$ws.Columns.Item(10).ColumnWidth = 13
$ws.Columns.Item(12).ColumnWidth = 16
$ws.Columns.Item(22).ColumnWidth = 12.8
$ws.Columns.Item(24).ColumnWidth = 10.65
$ws.Columns.Item(25).ColumnWidth = 10.65
$ws.Columns.Item(26).ColumnWidth = 10.65

# ---------------------------------------------------------------------------
# 5. Selection / view state.
# ---------------------------------------------------------------------------
$ws.Range("M16").Select()
